$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new row 32 with the next time log entry ("Milestone 5 work")
$ws.Range("A32").Value = 44631
$ws.Range("B32").Value = 44631.041666666664
$ws.Range("C32").Formula = "=B32-A32"
$ws.Range("D32").Formula = "=C32+D31"
$ws.Range("E32").Value = "Milestone 5 work"

# Copy formatting (number formats / styles) from the row above to keep things consistent
$ws.Range("A31:E31").Copy() | Out-Null
$ws.Range("A32:E32").PasteSpecial(-4122) | Out-Null # xlPasteFormats
$excel.CutCopyMode = 0

# Re-select the new active cell, matching where the cursor ends up after entry
$ws.Range("F32").Select() | Out-Null

$wb.Save()
